$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fmax(L)")
$rng = $ws.Range("K6:O6")
$rng.ClearFormats()
$b = $rng.Borders.Item(7)
$b.LineStyle = "Continuous"
$b.Weight = "Medium"
